$wb = $excel.ActiveWorkbook

# Update the conversion text on "Hoja1"
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 9.5 = 39267.74 pesos`n✅ 39267.74 pesos = 9.46 = 965.81 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# Update the rate values on "tasas"
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 105.277
$wsTasas.Range("O10").Value = 4133.99
$wsTasas.Range("N12").Value = 4150
$wsTasas.Range("O12").Value = 102.071

$wb.Save()
